$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks numeric but must be stored as TEXT
# (matches the existing shared-string "numbers-as-text" cells like the
# order-limit column and the sale-price column), while preserving the
# cell's original number format / style index.
function Set-TextValue($range, [string]$value, [string]$restoreFormat) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $restoreFormat
}

# ---------------------------------------------------------------------
# The report table (rows 7-19) is re-rendered with a new medicine,
# "DEPAKINE CHRONO 500MG 30 SCORED PROLONGED REL. F.C. TAB.", inserted
# in alphabetical order between "CYMBATEX ..." (row 11) and the old
# "DOLIPRANE ..." row (row 12). Every row from the old row 12 onward
# shifts down by one logical entry. Rather than performing a physical
# row-insert (which would also shift row heights/styles), we overwrite
# the cell contents of rows 12-20 directly, matching the target layout.
# ---------------------------------------------------------------------

# A-column item numbers 6..14, C/H/L/N/P/Q values, in their final order
$rows = @(
    @{ Row=12; Num=6;  Name="DEPAKINE CHRONO 500MG 30 SCORED PROLONGED REL. F.C. TAB."; H="1:0";  L="1"; N="144.00"; P="144.0000"; Q="1:0" },
    @{ Row=13; Num=7;  Name="DOLIPRANE 1 GM 15 TABS.";                                  H="8:3";  L="1"; N="48.00";  P="15.8400";  Q="0:1" },
    @{ Row=14; Num=8;  Name="EZACARD 75MG 30 ENTERIC COATED TAB.";                      H="1:1";  L="1"; N="51.00";  P="51.0000";  Q="1:0" },
    @{ Row=15; Num=9;  Name="FLUMOX 500MG 16 CAPS";                                     H="1:0";  L="1"; N="71.00";  P="35.5000";  Q="0:1" },
    @{ Row=16; Num=10; Name="LARYPRO 20 LOZENGES";                                      H="0:1";  L="1"; N="44.00";  P="22.0000";  Q="0:1" },
    @{ Row=17; Num=11; Name="LASILACTONE 50/20MG 30 TAB.";                              H="0:2";  L="1"; N="126.00"; P="126.0000"; Q="1:0" },
    @{ Row=18; Num=12; Name="RIVO 320MG 20*10 TABS";                                    H="1:9";  L="1"; N="141.00"; P="14.1000";  Q="0:2" },
    @{ Row=19; Num=13; Name="بادي سبلاش ايفا";                                          H="2:0";  L="0"; N="175.00"; P="350.0000"; Q="2:0" },
    @{ Row=20; Num=14; Name="صابون ديتول العنايه بالبشره";                              H="13:0"; L="0"; N="45.00";  P="45.0000";  Q="1:0" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Num                      # A - item number (numeric)
    $ws.Cells.Item($row, 3).Value = $r.Name                     # C - item name (text)
    $ws.Cells.Item($row, 8).Value = $r.H                        # H - balance "x:y" (text, non-numeric)
    Set-TextValue $ws.Cells.Item($row, 12) $r.L "#,##0.##;""[""#,##0.##""]"";0"   # L - order limit (numeric-looking text)
    $ws.Cells.Item($row, 14).Value = $r.N                       # N - price (text, format already text)
    Set-TextValue $ws.Cells.Item($row, 16) $r.P "0.00"          # P - sale price (numeric-looking text)
    $ws.Cells.Item($row, 17).Value = $r.Q                       # Q - transactions "x:y" (text, non-numeric)
}

# Row 20 is a brand-new table row; give it the same row height, merges,
# and per-cell styles as the other item rows.
$ws.Rows.Item(20).RowHeight = 24.75
$ws.Range("A20:B20").Merge()
$ws.Range("C20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()
$ws.Range("N20:O20").Merge()
$ws.Range("A20:Q20").Style = $ws.Range("A19:Q19").Style

# ---------------------------------------------------------------------
# New totals row (was row 20) moves to row 21, with the updated sum.
# ---------------------------------------------------------------------
$ws.Range("A21:Q21").Style = $ws.Range("A20:Q20").Style
$ws.Rows.Item(21).RowHeight = 25.5
$ws.Cells.Item(21, 16).Value = 1188.8599999999999
$ws.Cells.Item(21, 16).Style = $ws.Cells.Item(20, 16).Style
$ws.Cells.Item(21, 17).Style = $ws.Cells.Item(20, 17).Style
$ws.Range("P21:Q21").Merge()

# ---------------------------------------------------------------------
# New footer row (was row 21) moves to row 22, with the refreshed
# generation timestamp.
# ---------------------------------------------------------------------
$ws.Range("A22:Q22").Style = $ws.Range("A21:Q21").Style
$ws.Rows.Item(22).RowHeight = 16.5
$ws.Cells.Item(22, 1).Value = "Thursday, 25 September, 2025 1:00 PM"
$ws.Cells.Item(22, 7).Value = "1/1"
$ws.Cells.Item(22, 11).Value = "developed by : Abdelaziz Talaat"
$ws.Range("A22:F22").Merge()
$ws.Range("G22:I22").Merge()
$ws.Range("K22:Q22").Merge()

Write-Host "edit applied"
